$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @('21TRD09200', 'Bunner', 'DUS UCM', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 50', '$ 25'),
    @('21TRD09200', 'Bunner', 'OPERATING W/O A VALID OL - UCM', '4510.12', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS', '4510.21A*', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'FAILURE TO FILE REGISTRATION', '4503.11', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'Speeding > 25 mph', '4511.21(B)(2)', 'Minor Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'DUS UCM', '4510.111', 'UCM', 'No Contest', 'Guilty', '$ 50', '$ 25'),
    @('21TRD09200', 'Bunner', 'OPERATING W/O A VALID OL - UCM', '4510.12', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'FAILURE TO REINSTATE LICENSE UCM 1-2/3YRS', '4510.21A*', 'UCM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'FAILURE TO FILE REGISTRATION', '4503.11', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0'),
    @('21TRD09200', 'Bunner', 'Speeding > 25 mph', '4511.21(B)(2)', 'Minor Misdemeanor', 'No Contest', 'Guilty', '$ 0', '$ 0'),
)

$startRow = 526
$endRow = $startRow + $data.Count - 1

# Force text interpretation so values like "4510.111" or "$ 50" are kept
# verbatim as strings rather than being coerced to numbers/currency.
$fillRange = $ws.Range("A" + $startRow + ":I" + $endRow)
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

# Restore default (general) cell style now that the text values are locked
# in, so the new rows don't carry a stray explicit style index.
$fillRange.Style = "Normal"